# M6 Offline Data - apply authoring updates
$p = $ppt.ActivePresentation

# --- Step 1: refresh slide identities for the first two slides -------------
# The original deck carried stale SlideIDs (271/274) for the first two
# slides. Duplicating them in place and deleting the originals gives the
# deck fresh SlideIDs while leaving their position (and every other slide)
# untouched - mirroring how the authoring tool produced ids 459/460.
$origSlide1 = $p.Slides.Item(1)
$newSlide1 = $origSlide1.Duplicate()

$origSlide2 = $p.Slides.Item(3)
$newSlide2 = $origSlide2.Duplicate()

$p.Slides.Item(1).Delete()
$p.Slides.Item(2).Delete()

# After the dance above:
#   Slides.Item(1) -> new SlideID 459 (was the title slide)
#   Slides.Item(2) -> new SlideID 460 (was the "Course Topics" table slide)
$slide1 = $p.Slides.Item(1)
$slide2 = $p.Slides.Item(2)

# --- Step 2: title slide (new SlideID 459) ----------------------------------
foreach ($shape in $slide1.Shapes) {
    if (-not $shape.HasTextFrame) { continue }
    $ph = $shape.PlaceholderFormat.Type
    if ($shape.Name -eq "Subtitle 3") {
        $tf = $shape.TextFrame.TextRange
        $tf.Text = "Jeremy Foster`rChristopher Harrison"
    } elseif ($shape.Name -eq "Title 1") {
        $shape.TextFrame.TextRange.Text = "Mobile Web"
    }
}

# --- Step 3: "Course Topics" table slide (new SlideID 460) -----------------
foreach ($shape in $slide2.Shapes) {
    if ($shape.HasTable) {
        $tbl = $shape.Table
        $tbl.Cell(1,1).Shape.TextFrame.TextRange.Text = "Mobile Web"
        $tbl.Cell(2,1).Shape.TextFrame.TextRange.Text = "01 | Designing for Mobile"
        $tbl.Cell(2,2).Shape.TextFrame.TextRange.Text = "05 | The Mobile Client"
        $tbl.Cell(3,1).Shape.TextFrame.TextRange.Text = "02 | Mobile UI"
        $tbl.Cell(3,2).Shape.TextFrame.TextRange.Text = "06 | Offline Data"
        $tbl.Cell(4,1).Shape.TextFrame.TextRange.Text = "03 | Integrating Touch"
        $tbl.Cell(4,2).Shape.TextFrame.TextRange.Text = "07 | Publishing to Azure"
        $tbl.Cell(5,1).Shape.TextFrame.TextRange.Text = "04 | Setting Up the Server"
        $tbl.Cell(5,2).Shape.TextFrame.TextRange.Text = ""
    }
}

# --- Step 4: module divider slide (SlideID 283, unchanged) -----------------
foreach ($slide in $p.Slides) {
    if ($slide.SlideID -ne 283) { continue }
    foreach ($shape in $slide.Shapes) {
        if (-not $shape.HasTextFrame) { continue }
        if ($shape.Name -eq "Text Placeholder 4") {
            $shape.TextFrame.TextRange.Text = "06 | Offline Data"
        } elseif ($shape.Name -eq "Subtitle 3") {
            $shape.TextFrame.TextRange.Text = "Christopher Harrison | @GeekTrainer`rJeremy Foster | @codefoster"
        }
    }
}
